# TST: Update tests for changes to Area and Problem
#
# Reproduces the workbook edit: the "Problems" table on the "Area" sheet
# grows from 6 rows (Problem 1-6) to 25 rows (Problem 1-25), the boulder
# assignment / grade data for the existing rows 29-34 is corrected, and the
# sheet's selection moves down to track the newly-added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix up the existing "Problem" rows (29-34): boulder name + grade data
#    was re-worked so that problems 1-6 all belong to "Boulder A".
# ---------------------------------------------------------------------
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 1

$ws.Range("D30").Value = 2

$ws.Range("B31").Value = "Boulder A"
$ws.Range("C31").Value = 0

$ws.Range("B32").Value = "Boulder A"
$ws.Range("C32").Value = 3

$ws.Range("B33").Value = "Boulder A"
$ws.Range("C33").Value = 3

$ws.Range("B34").Value = "Boulder A"
$ws.Range("C34").Value = 4
$ws.Range("D34").Value = 2

# ---------------------------------------------------------------------
# 2. Append the new "Problem" rows (35-53): Problem 7-25, covering
#    "Boulder B" (7-19) and "Boulder C" (20-25). Column E keeps reusing
#    "Sixth boulder", matching a fill-down from row 34.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 35; Problem = "Problem 7";  Boulder = "Boulder B"; C = 0;  D = 3 },
    @{ Row = 36; Problem = "Problem 8";  Boulder = "Boulder B"; C = 6;  D = 2 },
    @{ Row = 37; Problem = "Problem 9";  Boulder = "Boulder B"; C = 6;  D = 5 },
    @{ Row = 38; Problem = "Problem 10"; Boulder = "Boulder B"; C = 7;  D = 1 },
    @{ Row = 39; Problem = "Problem 11"; Boulder = "Boulder B"; C = 8;  D = 3 },
    @{ Row = 40; Problem = "Problem 12"; Boulder = "Boulder B"; C = 11; D = 3 },
    @{ Row = 41; Problem = "Problem 13"; Boulder = "Boulder B"; C = 1;  D = 2 },
    @{ Row = 42; Problem = "Problem 14"; Boulder = "Boulder B"; C = 13; D = 3 },
    @{ Row = 43; Problem = "Problem 15"; Boulder = "Boulder B"; C = 5;  D = 4 },
    @{ Row = 44; Problem = "Problem 16"; Boulder = "Boulder B"; C = 6;  D = 5 },
    @{ Row = 45; Problem = "Problem 17"; Boulder = "Boulder B"; C = 7;  D = 4 },
    @{ Row = 46; Problem = "Problem 18"; Boulder = "Boulder B"; C = 0;  D = 4 },
    @{ Row = 47; Problem = "Problem 19"; Boulder = "Boulder B"; C = 9;  D = 3 },
    @{ Row = 48; Problem = "Problem 20"; Boulder = "Boulder C"; C = 7;  D = 2 },
    @{ Row = 49; Problem = "Problem 21"; Boulder = "Boulder C"; C = 7;  D = 2 },
    @{ Row = 50; Problem = "Problem 22"; Boulder = "Boulder C"; C = 4;  D = 3 },
    @{ Row = 51; Problem = "Problem 23"; Boulder = "Boulder C"; C = 6;  D = 3 },
    @{ Row = 52; Problem = "Problem 24"; Boulder = "Boulder C"; C = 2;  D = 4 },
    @{ Row = 53; Problem = "Problem 25"; Boulder = "Boulder C"; C = 6;  D = 5 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $r.Problem
    $ws.Range("B$rowNum").Value = $r.Boulder
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = "Sixth boulder"
}

# ---------------------------------------------------------------------
# 3. Move the view / selection to track the newly added rows, like the
#    author scrolled down to B36 after adding the data.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
try { $win.ScrollRow = 27 } catch { }
try { $win.ScrollColumn = 1 } catch { }

$ws.Range("B36").Select()
